# "All changes as of 7-27"
# Inserts 10 new codebook rows describing ACS "means of transportation to
# work by earnings" / "worked from home" variables into the Codebook sheet,
# right after the existing "Pop_Grad" (B06009_006E) entry and before the
# "Less_High_pct" block, pushing that block (and the blank rows after it)
# down by 10 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. The old row 113 (Pop_Grad) carried a medium bottom border marking the
#    end of the "education" section box. That boundary is moving down to
#    the new last row of the inserted block (new row 123), so strip it off
#    row 113 first.
# ---------------------------------------------------------------------
$ws.Range("A113:J113").Borders.Item(9).LineStyle = -4142

# ---------------------------------------------------------------------
# 2. Insert 10 blank rows right before the old row 114 ("Less_High_pct"),
#    shifting it (and everything below, through the trailing blank rows)
#    down to row 124 onward.
# ---------------------------------------------------------------------
$ws.Range("A114:J123").EntireRow.Insert()

# ---------------------------------------------------------------------
# 3. Populate the 9 new "worked from home by earnings" rows (114-122) plus
#    the new "Median_Income_Worked_home" row (123), which becomes the new
#    end-of-section boundary row.
# ---------------------------------------------------------------------
$rows = @(
  @{ r=114; a="Pop_Worked_Home";          c="B08119_055E"; j="MEANS OF TRANSPORTATION TO WORK BY WORKERS' EARNINGS IN THE PAST 12 MONTHS (IN 2021 INFLATION-ADJUSTED DOLLARS)"; h="sum" },
  @{ r=115; a="Pop_Worked_Home_1-10k";     c="B08119_056E"; j="MEANS OF TRANSPORTATION TO WORK BY WORKERS' EARNINGS IN THE PAST 12 MONTHS (IN 2021 INFLATION-ADJUSTED DOLLARS)"; h="sum" },
  @{ r=116; a="Pop_Worked_Home_10k-15k";   c="B08119_057E"; j="MEANS OF TRANSPORTATION TO WORK BY WORKERS' EARNINGS IN THE PAST 12 MONTHS (IN 2021 INFLATION-ADJUSTED DOLLARS)"; h="sum" },
  @{ r=117; a="Pop_Worked_Home_15k-25k";   c="B08119_058E"; j="MEANS OF TRANSPORTATION TO WORK BY WORKERS' EARNINGS IN THE PAST 12 MONTHS (IN 2021 INFLATION-ADJUSTED DOLLARS)"; h="sum" },
  @{ r=118; a="Pop_Worked_Home_25k-35k";   c="B08119_059E"; j="MEANS OF TRANSPORTATION TO WORK BY WORKERS' EARNINGS IN THE PAST 12 MONTHS (IN 2021 INFLATION-ADJUSTED DOLLARS)"; h="sum" },
  @{ r=119; a="Pop_Worked_Home_35k-50k";   c="B08119_060E"; j="MEANS OF TRANSPORTATION TO WORK BY WORKERS' EARNINGS IN THE PAST 12 MONTHS (IN 2021 INFLATION-ADJUSTED DOLLARS)"; h="sum" },
  @{ r=120; a="Pop_Worked_Home_50k-65k";   c="B08119_061E"; j="MEANS OF TRANSPORTATION TO WORK BY WORKERS' EARNINGS IN THE PAST 12 MONTHS (IN 2021 INFLATION-ADJUSTED DOLLARS)"; h="sum" },
  @{ r=121; a="Pop_Worked_Home_65k-75k";   c="B08119_062E"; j="MEANS OF TRANSPORTATION TO WORK BY WORKERS' EARNINGS IN THE PAST 12 MONTHS (IN 2021 INFLATION-ADJUSTED DOLLARS)"; h="sum" },
  @{ r=122; a="Pop_Worked_Home_75k-more";  c="B08119_063E"; j="MEANS OF TRANSPORTATION TO WORK BY WORKERS' EARNINGS IN THE PAST 12 MONTHS (IN 2021 INFLATION-ADJUSTED DOLLARS)"; h="sum" },
  @{ r=123; a="Median_Income_Worked_home"; c="B08121_007E"; j="MEDIAN EARNINGS IN THE PAST 12 MONTHS (IN 2021 INFLATION-ADJUSTED DOLLARS) BY MEANS OF TRANSPORTATION TO WORK"; h="mean" }
)

foreach ($row in $rows) {
  $r = $row.r
  $ws.Cells.Item($r, 1).Value = $row.a          # A: column name
  $ws.Cells.Item($r, 2).Value = "ACS"            # B: data source
  $ws.Cells.Item($r, 3).Value = $row.c           # C: var code
  $ws.Cells.Item($r, 4).Value = "2010, 2015, 2019, 2021"  # D: year(s)
  $ws.Cells.Item($r, 8).Value = $row.h           # H: agg type
  $ws.Cells.Item($r, 10).Value = $row.j          # J: description
}

# ---------------------------------------------------------------------
# 4. Apply the section-boundary bottom border to the new closing row 123
#    (same formatting the old row 113 used to carry).
# ---------------------------------------------------------------------
$border = $ws.Range("A123:J123").Borders.Item(9)
$border.LineStyle = 1
$border.Weight = -4138
$border.Color = 0
